$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 4506
$ws.Range("J86").Value = 5324.6665
$ws.Range("L86").Value = 5324.6665
$ws.Range("N86").Value = -7570.6665
$ws.Range("H89").Value = 4506
$ws.Range("J89").Value = 5324.6665
$ws.Range("L89").Value = 26623.3325
$ws.Range("N89").Value = -37855.3325
$ws.Range("H132").Value = 1589.0834
$ws.Range("I132").Value = 1381.0938
$ws.Range("K132").Value = 4143.2814
$ws.Range("M132").Value = -1613.2814
$ws.Range("H138").Value = 3238.6191
$ws.Range("I138").Value = 2427.625
$ws.Range("J138").Value = 4319.9443
$ws.Range("K138").Value = 7282.875
$ws.Range("L138").Value = 12959.8329
$ws.Range("M138").Value = -2142.875
$ws.Range("N138").Value = -23239.8329
$ws.Range("H141").Value = 744.9231
$ws.Range("I141").Value = 744.9231
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 2234.7693
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = 2945.2307
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3021.7437
$ws.Range("I32").Value = 2395.1177
$ws.Range("J32").Value = 7282.8
$ws.Range("K32").Value = 2395.1177
$ws.Range("L32").Value = 7282.8
$ws.Range("M32").Value = -2108.1177
$ws.Range("N32").Value = -7856.8
$ws.Range("H61").Value = 8274.966
$ws.Range("I61").Value = 6533.88
$ws.Range("J61").Value = 19156.75
$ws.Range("K61").Value = 6533.88
$ws.Range("L61").Value = 19156.75
$ws.Range("M61").Value = -6321.88
$ws.Range("N61").Value = -19580.75
$ws.Range("H74").Value = 25643462
$ws.Range("I74").Value = 27780250
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 27780250
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -27779376
$ws.Range("N74").Value = -3748
$ws.Range("H77").Value = 25643462
$ws.Range("I77").Value = 27780250
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 138901250
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -138896882
$ws.Range("N77").Value = -18736
$ws.Range("H136").Value = 8274.966
$ws.Range("I136").Value = 6533.88
$ws.Range("J136").Value = 19156.75
$ws.Range("K136").Value = 19601.64
$ws.Range("L136").Value = 57470.25
$ws.Range("M136").Value = -17051.64
$ws.Range("N136").Value = -62570.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1898.32
$ws.Range("I20").Value = 1410.5294
$ws.Range("K20").Value = 1410.5294
$ws.Range("M20").Value = -1163.5294
$ws.Range("H64").Value = 1950.8572
$ws.Range("I64").Value = 1941.2
$ws.Range("J64").Value = 1975
$ws.Range("K64").Value = 1941.2
$ws.Range("L64").Value = 1975
$ws.Range("M64").Value = -1716.2
$ws.Range("N64").Value = -2425
$ws.Range("H67").Value = 1950.8572
$ws.Range("I67").Value = 1941.2
$ws.Range("J67").Value = 1975
$ws.Range("K67").Value = 1941.2
$ws.Range("L67").Value = 1975
$ws.Range("M67").Value = -1161.2
$ws.Range("N67").Value = -3535
$ws.Range("H86").Value = 2867.2222
$ws.Range("I86").Value = 2901.875
$ws.Range("J86").Value = 2590
$ws.Range("K86").Value = 2901.875
$ws.Range("L86").Value = 2590
$ws.Range("M86").Value = -1778.875
$ws.Range("N86").Value = -4836
$ws.Range("H89").Value = 2867.2222
$ws.Range("I89").Value = 2901.875
$ws.Range("J89").Value = 2590
$ws.Range("K89").Value = 14509.375
$ws.Range("L89").Value = 12950
$ws.Range("M89").Value = -8893.375
$ws.Range("N89").Value = -24182
$ws.Range("H135").Value = 51999.625
$ws.Range("J135").Value = 51999.625
$ws.Range("L135").Value = 51999.625
$ws.Range("N135").Value = -62139.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 32821.688
$ws.Range("I31").Value = 2832.913
$ws.Range("J31").Value = 90300.164
$ws.Range("K31").Value = 2832.913
$ws.Range("L31").Value = 90300.164
$ws.Range("M31").Value = -2537.913
$ws.Range("N31").Value = -90890.164
$ws.Range("H34").Value = 32821.688
$ws.Range("I34").Value = 2832.913
$ws.Range("J34").Value = 90300.164
$ws.Range("K34").Value = 2832.913
$ws.Range("L34").Value = 90300.164
$ws.Range("M34").Value = -2630.913
$ws.Range("N34").Value = -90704.164
$ws.Range("H58").Value = 3385.12
$ws.Range("I58").Value = 1287.1052
$ws.Range("K58").Value = 1287.1052
$ws.Range("M58").Value = -1084.1052
$ws.Range("H132").Value = 13598
$ws.Range("I132").Value = 15659.333
$ws.Range("K132").Value = 46977.999
$ws.Range("M132").Value = -44447.999
$ws.Range("H134").Value = 3189.55
$ws.Range("I134").Value = 2003.3846
$ws.Range("J134").Value = 5392.4287
$ws.Range("K134").Value = 6010.1538
$ws.Range("L134").Value = 16177.2861
$ws.Range("M134").Value = -3475.1538
$ws.Range("N134").Value = -21247.2861
$ws.Range("H136").Value = 3385.12
$ws.Range("I136").Value = 1287.1052
$ws.Range("K136").Value = 3861.3156
$ws.Range("M136").Value = -1311.3156

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3920
$ws.Range("I81").Value = 1378.5714
$ws.Range("K81").Value = 4135.7142
$ws.Range("M81").Value = -3012.7142
$ws.Range("H84").Value = 3920
$ws.Range("I84").Value = 1378.5714
$ws.Range("K84").Value = 12407.1426
$ws.Range("M84").Value = -6791.142600000001
$ws.Range("H133").Value = 5409.5713
$ws.Range("I133").Value = 3560.2
$ws.Range("J133").Value = 10033
$ws.Range("K133").Value = 10680.6
$ws.Range("L133").Value = 30099
$ws.Range("M133").Value = -5620.599999999999
$ws.Range("N133").Value = -40219
$ws.Range("H140").Value = 2083.36
$ws.Range("I140").Value = 1367.9546
$ws.Range("K140").Value = 4103.8638
$ws.Range("M140").Value = 1076.1362
$ws.Range("H141").Value = 4848.091
$ws.Range("I141").Value = 3429.6
$ws.Range("K141").Value = 10288.8
$ws.Range("M141").Value = -5108.799999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 10000
$ws.Range("J34").Value = 10000
$ws.Range("L34").Value = 10000
$ws.Range("N34").Value = -10536
$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()
$ws.Range("H64").Value = 19499.5
$ws.Range("J64").Value = 19499.5
$ws.Range("L64").Value = 19499.5
$ws.Range("N64").Value = -19995.5
$ws.Range("H67").Value = 19499.5
$ws.Range("J67").Value = 19499.5
$ws.Range("L67").Value = 19499.5
$ws.Range("N67").Value = -21215.5
$ws.Range("H70").Value = 7797.636
$ws.Range("I70").Value = 5809.5625
$ws.Range("J70").Value = 13099.167
$ws.Range("K70").Value = 5809.5625
$ws.Range("L70").Value = 13099.167
$ws.Range("M70").Value = -5539.5625
$ws.Range("N70").Value = -13639.167
$ws.Range("H73").Value = 7797.636
$ws.Range("I73").Value = 5809.5625
$ws.Range("J73").Value = 13099.167
$ws.Range("K73").Value = 5809.5625
$ws.Range("L73").Value = 13099.167
$ws.Range("M73").Value = -4873.5625
$ws.Range("N73").Value = -14971.167
$ws.Range("H76").Value = 10000
$ws.Range("J76").Value = 10000
$ws.Range("L76").Value = 10000
$ws.Range("N76").Value = -10630
$ws.Range("H79").Value = 10000
$ws.Range("J79").Value = 10000
$ws.Range("L79").Value = 10000
$ws.Range("N79").Value = -12184
$ws.Range("H132").Value = 4139.5713
$ws.Range("I132").Value = 2849.75
$ws.Range("J132").Value = 11878.5
$ws.Range("K132").Value = 8549.25
$ws.Range("L132").Value = 35635.5
$ws.Range("M132").Value = -6019.25
$ws.Range("N132").Value = -40695.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 8844.556
$ws.Range("J22").Value = 12383.5
$ws.Range("L22").Value = 12383.5
$ws.Range("N22").Value = -12973.5
$ws.Range("H27").Value = 8844.556
$ws.Range("J27").Value = 12383.5
$ws.Range("L27").Value = 12383.5
$ws.Range("N27").Value = -12597.5
$ws.Range("H46").Value = 2981.1875
$ws.Range("I46").Value = 767.3333
$ws.Range("J46").Value = 3492.077
$ws.Range("K46").Value = 767.3333
$ws.Range("L46").Value = 3492.077
$ws.Range("M46").Value = -579.3333
$ws.Range("N46").Value = -3868.077
$ws.Range("H55").Value = 1298.7407
$ws.Range("I55").Value = 278.4
$ws.Range("J55").Value = 2574.1667
$ws.Range("K55").Value = 278.4
$ws.Range("L55").Value = 2574.1667
$ws.Range("M55").Value = -105.4
$ws.Range("N55").Value = -2920.1667
$ws.Range("H74").Value = 47497.5
$ws.Range("I74").Value = 39995
$ws.Range("K74").Value = 39995
$ws.Range("M74").Value = -38997
$ws.Range("H77").Value = 47497.5
$ws.Range("I77").Value = 39995
$ws.Range("K77").Value = 119985
$ws.Range("M77").Value = -114993
$ws.Range("H132").Value = 5509.56
$ws.Range("I132").Value = 2735.2222
$ws.Range("J132").Value = 12643.571
$ws.Range("K132").Value = 8205.6666
$ws.Range("L132").Value = 37930.713
$ws.Range("M132").Value = -5675.6666
$ws.Range("N132").Value = -42990.713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 5965.353
$ws.Range("I132").Value = 5284.4194
$ws.Range("K132").Value = 15853.2582
$ws.Range("M132").Value = -13323.2582

Write-Output "Applied all changes"